$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.25
$ws.Range("G2").Value = 3.3
$ws.Range("H2").Value = 2.36
$ws.Range("I2").Value = 2.42
$ws.Range("J2").Value = 3.55
$ws.Range("L2").Value = 1.43
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.9
$ws.Range("P2").Value = 1.95
$ws.Range("R2").Value = 1.37
$ws.Range("S2").Value = 3.55
$ws.Range("T2").Value = 1.79
$ws.Range("U2").Value = 2.18
$ws.Range("V2").Value = 1.71
$ws.Range("W2").Value = 1.43
$ws.Range("X2").Value = 14
$ws.Range("Y2").Value = 10.5
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 13
$ws.Range("AC2").Value = 7.8
$ws.Range("AE2").Value = 25
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 14
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 48
$ws.Range("AJ2").Value = 60
$ws.Range("AK2").Value = 38
$ws.Range("AL2").Value = 48
$ws.Range("AN2").Value = 36
$ws.Range("AO2").Value = 21
